# Insert a new row above row 176, shifting existing rows 176-183 down to 177-184.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(176).Insert()

# Populate the newly inserted row 176 with values (copy constant columns from row 177,
# which now holds the data that used to be in row 176 before the shift).
$ws.Range("A176").Value = 7
$ws.Range("B176").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C176").Value = "Ñuble"
$ws.Range("D176").Value = 44509
$ws.Range("D176").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E176").Value = 16
$ws.Range("F176").Value = 100112008
$ws.Range("G176").Value = "Coliflor"
$ws.Range("H176").Value = "Sin especificar"
$ws.Range("I176").Value = "Primera"
$ws.Range("J176").Value = 300
$ws.Range("K176").Value = 700
$ws.Range("L176").Value = 750
$ws.Range("M176").Value = 725
$ws.Range("N176").Value = '$/unidad'
$ws.Range("O176").Value = "Región del Maule"
$ws.Range("P176").Value = 725
$ws.Range("Q176").Value = 1
$ws.Range("R176").Value = "Hortaliza"
